$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 349
$ws.Range("I19").Value = 300
$ws.Range("K19").Value = 300
$ws.Range("M19").Value = -125
# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 6247.654
$ws.Range("I64").Value = 3537.6667
$ws.Range("J64").Value = 7682.353
$ws.Range("K64").Value = 3537.6667
$ws.Range("L64").Value = 7682.353
$ws.Range("M64").Value = -3289.6667
$ws.Range("N64").Value = -8178.353
# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 6247.654
$ws.Range("I67").Value = 3537.6667
$ws.Range("J67").Value = 7682.353
$ws.Range("K67").Value = 3537.6667
$ws.Range("L67").Value = 7682.353
$ws.Range("M67").Value = -2679.6667
$ws.Range("N67").Value = -9398.352999999999
# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 3778.0557
$ws.Range("I86").Value = 3850.3
$ws.Range("J86").Value = 3687.75
$ws.Range("K86").Value = 3850.3
$ws.Range("L86").Value = 3687.75
$ws.Range("M86").Value = -2727.3
$ws.Range("N86").Value = -5933.75
# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 3778.0557
$ws.Range("I89").Value = 3850.3
$ws.Range("J89").Value = 3687.75
$ws.Range("K89").Value = 19251.5
$ws.Range("L89").Value = 18438.75
$ws.Range("M89").Value = -13635.5
$ws.Range("N89").Value = -29670.75
# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 1806.1666
$ws.Range("J129").Value = 3750
$ws.Range("L129").Value = 11250
$ws.Range("N129").Value = -21250
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2049.8333
$ws.Range("J138").Value = 2852.3333
$ws.Range("L138").Value = 8556.999899999999
$ws.Range("N138").Value = -18836.9999
# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 5168.375
$ws.Range("I141").Value = 2520.8948
$ws.Range("J141").Value = 9037.77
$ws.Range("K141").Value = 7562.6844
$ws.Range("L141").Value = 27113.31
$ws.Range("M141").Value = -2382.6844
$ws.Range("N141").Value = -37473.31

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 34 (Leve Item ID 2753)
$ws.Range("H34").Value = 213332.67
$ws.Range("J34").Value = 213332.67
$ws.Range("L34").Value = 213332.67
$ws.Range("N34").Value = -213874.67
# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 8133.3335
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 8133.3335
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 4045.5
$ws.Range("J74").Value = 4240.6665
$ws.Range("L74").Value = 4240.6665
$ws.Range("N74").Value = -5988.6665
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 4045.5
$ws.Range("J77").Value = 4240.6665
$ws.Range("L77").Value = 21203.3325
$ws.Range("N77").Value = -29939.3325
# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 1410.5
$ws.Range("I88").Value = 6
$ws.Range("J88").Value = 1566.5555
$ws.Range("K88").Value = 6
$ws.Range("L88").Value = 1566.5555
$ws.Range("M88").Value = 400
$ws.Range("N88").Value = -2378.5555
# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 1410.5
$ws.Range("I91").Value = 6
$ws.Range("J91").Value = 1566.5555
$ws.Range("K91").Value = 6
$ws.Range("L91").Value = 1566.5555
$ws.Range("M91").Value = 1398
$ws.Range("N91").Value = -4374.5555
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2948.9
$ws.Range("I122").Value = 2914.8333
$ws.Range("K122").Value = 8744.499899999999
$ws.Range("M122").Value = -6294.499899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 18520354
$ws.Range("I86").Value = 1794.4546
$ws.Range("J86").Value = 100002020
$ws.Range("K86").Value = 1794.4546
$ws.Range("L86").Value = 100002020
$ws.Range("M86").Value = -671.4546
$ws.Range("N86").Value = -100004266
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 18520354
$ws.Range("I89").Value = 1794.4546
$ws.Range("J89").Value = 100002020
$ws.Range("K89").Value = 8972.273000000001
$ws.Range("L89").Value = 500010100
$ws.Range("M89").Value = -3356.273000000001
$ws.Range("N89").Value = -500021332
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 3738.6924
$ws.Range("I94").Value = 3858.4285
$ws.Range("K94").Value = 3858.4285
$ws.Range("M94").Value = -3407.4285
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 4259.95
$ws.Range("J99").Value = 5245.5
$ws.Range("L99").Value = 5245.5
$ws.Range("N99").Value = -8241.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4886.6665
$ws.Range("I31").Value = 2081.5334
$ws.Range("J31").Value = 11899.5
$ws.Range("K31").Value = 2081.5334
$ws.Range("L31").Value = 11899.5
$ws.Range("M31").Value = -1786.5334
$ws.Range("N31").Value = -12489.5
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4886.6665
$ws.Range("I34").Value = 2081.5334
$ws.Range("J34").Value = 11899.5
$ws.Range("K34").Value = 2081.5334
$ws.Range("L34").Value = 11899.5
$ws.Range("M34").Value = -1879.5334
$ws.Range("N34").Value = -12303.5
# Row 74 (Leve Item ID 10636)
$ws.Range("H74").Value = 41666
$ws.Range("J74").Value = 41666
$ws.Range("L74").Value = 41666
$ws.Range("N74").Value = -43414
# Row 77 (Leve Item ID 10636)
$ws.Range("H77").Value = 41666
$ws.Range("J77").Value = 41666
$ws.Range("L77").Value = 124998
$ws.Range("N77").Value = -133734
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 3752.111
$ws.Range("I132").Value = 3596.125
$ws.Range("K132").Value = 10788.375
$ws.Range("M132").Value = -8258.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 64 (Leve Item ID 12861)
$ws.Range("H64").Value = 3504.5
$ws.Range("J64").Value = 3504.5
$ws.Range("L64").Value = 10513.5
$ws.Range("N64").Value = -11053.5
# Row 67 (Leve Item ID 12861)
$ws.Range("H67").Value = 3504.5
$ws.Range("J67").Value = 3504.5
$ws.Range("L67").Value = 10513.5
$ws.Range("N67").Value = -12385.5
# Row 86 (Leve Item ID 12892)
$ws.Range("H86").Value = 378
$ws.Range("I86").Value = 296.66666
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 889.9999799999999
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = 296.0000200000001
$ws.Range("N86").Value = -3872
# Row 89 (Leve Item ID 12892)
$ws.Range("H89").Value = 378
$ws.Range("I89").Value = 296.66666
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 2669.99994
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 3258.00006
$ws.Range("N89").Value = -16356
# Row 94 (Leve Item ID 19811)
$ws.Range("H94").Value = 8633.223
$ws.Range("J94").Value = 9337.375
$ws.Range("L94").Value = 28012.125
$ws.Range("N94").Value = -29364.125
# Row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 671.3570999999999
$ws.Range("I121").Value = 536.36365
$ws.Range("J121").Value = 1166.3334
$ws.Range("K121").Value = 1609.09095
$ws.Range("L121").Value = 3499.0002
$ws.Range("M121").Value = -299.09095
$ws.Range("N121").Value = -6119.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 104 (Leve Item ID 18666)
$ws.Range("H104").Value = 198999
$ws.Range("J104").Value = 198999
$ws.Range("L104").Value = 198999
$ws.Range("N104").Value = -205987
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 9803.879999999999
$ws.Range("I122").Value = 9803.879999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 29411.64
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -26961.64
$ws.Range("N122").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 4099.8
$ws.Range("I40").Value = 2592.1428
$ws.Range("J40").Value = 6361.2856
$ws.Range("K40").Value = 2592.1428
$ws.Range("L40").Value = 6361.2856
$ws.Range("M40").Value = -2456.1428
$ws.Range("N40").Value = -6633.2856
# Row 62 (Leve Item ID 10740)
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
# Row 65 (Leve Item ID 10740)
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 7399.7334
$ws.Range("I93").Value = 6299.2
$ws.Range("K93").Value = 6299.2
$ws.Range("M93").Value = -5051.2
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 3819.2
$ws.Range("I122").Value = 3819.2
$ws.Range("K122").Value = 11457.6
$ws.Range("M122").Value = -9007.599999999999
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3662.6287
$ws.Range("I132").Value = 3354.8965
$ws.Range("K132").Value = 10064.6895
$ws.Range("M132").Value = -7534.6895

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 41 (Leve Item ID 21725)
$ws.Range("H41").Value = 13598.4
$ws.Range("I41").Value = 20000
$ws.Range("J41").Value = 11998
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 11998
$ws.Range("M41").Value = -19610
$ws.Range("N41").Value = -12778
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 4848.294
$ws.Range("I122").Value = 1046.7778
$ws.Range("J122").Value = 9125
$ws.Range("K122").Value = 3140.3334
$ws.Range("L122").Value = 27375
$ws.Range("M122").Value = -690.3334000000004
$ws.Range("N122").Value = -32275
